# Generate Report for Handback
# Adds a new handback record (e089686e-1a2e-4984-bed7-c10b48441f07) as row 4
# on each of the three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId = "e089686e-1a2e-4984-bed7-c10b48441f07"
$mdName = "$fileId.md"
$zhXlf = "$fileId.34bc86967529de48652e605f86d015614d3839b6.zh-cn.xlf"
$deXlf = "$fileId.34bc86967529de48652e605f86d015614d3839b6.de-de.xlf"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $tblOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1, 1).Value = $mdName
$rngOverview.Cells.Item(1, 3).Value = ".md"
$rngOverview.Cells.Item(1, 5).Value = $statusText
$rngOverview.Cells.Item(1, 6).Value = $statusText
$rngOverview.Cells.Item(1, 7).Value = "2016-08-17 14:44:58"
$rngOverview.Cells.Item(1, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($rngOverview.Cells.Item(1, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0db03b7f0cad8a8f4196144d388be4e47a544c3a/e2e/$mdName", "", "", "e2e\$mdName")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $tblZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$rngZhCn.Cells.Item(1, 2).Value = ".md"
$rngZhCn.Cells.Item(1, 3).Value = $statusText
$rngZhCn.Cells.Item(1, 4).Value = "e2e"
$rngZhCn.Cells.Item(1, 5).Value = "ht"
$rngZhCn.Cells.Item(1, 6).Value = "True"
$rngZhCn.Cells.Item(1, 7).Value = $zhXlf
$rngZhCn.Cells.Item(1, 8).Value = "2016-08-17 14:44:53"
$rngZhCn.Cells.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngZhCn.Cells.Item(1, 10).Value = $zhXlf
$rngZhCn.Cells.Item(1, 11).Value = "2016-08-17 14:45:30"
$rngZhCn.Cells.Item(1, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngZhCn.Cells.Item(1, 13).Value = "True"
$rngZhCn.Cells.Item(1, 15).Value = "False"

$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0db03b7f0cad8a8f4196144d388be4e47a544c3a/e2e/$mdName", "", "", $mdName)
$wsZhCn.Hyperlinks.Add($rngZhCn.Cells.Item(1, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/22994a140d6b2eda91d5ec82eed685eab0acd942/e2e/$mdName", "", "", $mdName)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $tblDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$rngDeDe.Cells.Item(1, 2).Value = ".md"
$rngDeDe.Cells.Item(1, 3).Value = $statusText
$rngDeDe.Cells.Item(1, 4).Value = "e2e"
$rngDeDe.Cells.Item(1, 5).Value = "ht"
$rngDeDe.Cells.Item(1, 6).Value = "True"
$rngDeDe.Cells.Item(1, 7).Value = $deXlf
$rngDeDe.Cells.Item(1, 8).Value = "2016-08-17 14:44:58"
$rngDeDe.Cells.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngDeDe.Cells.Item(1, 10).Value = $deXlf
$rngDeDe.Cells.Item(1, 11).Value = "2016-08-17 14:45:37"
$rngDeDe.Cells.Item(1, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngDeDe.Cells.Item(1, 13).Value = "True"
$rngDeDe.Cells.Item(1, 15).Value = "False"

$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0db03b7f0cad8a8f4196144d388be4e47a544c3a/e2e/$mdName", "", "", $mdName)
$wsDeDe.Hyperlinks.Add($rngDeDe.Cells.Item(1, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f0303bb4519540cad9a80f3e29a0a1fc9d36ba3c/e2e/$mdName", "", "", $mdName)
